# "Updates in data followers"
# Adds a new column C ("ok") for the last block of rows (80-88, the
# "Pierre" group) and touches C65 (empty) so the underlying row-block
# metadata for rows 65-79 is refreshed (spans 1:2 -> 1:3) the same way the
# source workbook's XML shows, without actually adding a visible value
# there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch C65 to force the row-block (65-80) span metadata to widen to
# include column C, then clear it back out so no cell value is left in
# that row - matches rows 65-79 in the diff (spans changed, no new cell).
$ws.Range("C65").Value = "ok"
$ws.Range("C65").Value = ""

# Rows 80-88 (the "Pierre" / last group) each get a new column C value.
for ($r = 80; $r -le 88; $r++) {
    $ws.Range("C$r").Value = "ok"
}

# Restore the view's selection to match the edited workbook.
$ws.Range("E85").Select() | Out-Null
